# Automatic update of files.
# Applies row-level corrections to the "Artfynd" sheet:
#  - Increment Taxonsorteringsordning (column B) counters for several rows.
#  - Re-sort / correct a block of rows (10-14 and 23-25) so that each
#    observation's full record (Id, B, red-list status, taxon info,
#    coordinates, external id, times, observers, count/unit) lines up
#    with the correct row position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 91809
$ws.Range("B6").Value = 91809
$ws.Range("B7").Value = 79244
$ws.Range("B8").Value = 91809
$ws.Range("B9").Value = 92180
$ws.Range("B15").Value = 91809
$ws.Range("B16").Value = 91809
$ws.Range("B19").Value = 92107
$ws.Range("B20").Value = 91809
$ws.Range("B21").Value = 91809
$ws.Range("B22").Value = 91809
$ws.Range("B27").Value = 92268
$ws.Range("B28").Value = 92107
# Row 10 <- source row 12 (id 131106321)
$ws.Range("A10").Value = 131106321
$ws.Range("B10").Value = 92022
$ws.Range("D10").Value = "VU"
$ws.Range("E10").Value = 48
$ws.Range("F10").Value = "Lappticka"
$ws.Range("G10").Value = "Amylocystis lapponica"
$ws.Range("H10").Value = "(Romell) Bondartsev & Singer"
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = "mycel"
$ws.Range("Q10").Value = 601579
$ws.Range("R10").Value = 6992698
$ws.Range("X10").Value = "2025_0862"
$ws.Range("Z10").Value = "12:56"
$ws.Range("AB10").Value = "12:56"
$ws.Range("AX10").Value = "Alexander Hoffmann"

# Row 11 <- source row 10 (id 131106315)
$ws.Range("A11").Value = 131106315
$ws.Range("B11").Value = 92107
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 658
$ws.Range("F11").Value = "Rosenticka"
$ws.Range("G11").Value = "Fomitopsis rosea"
$ws.Range("H11").Value = "(Alb. & Schwein.:Fr.) P.Karst."
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "1"
$ws.Range("J11").Value = "mycel"
$ws.Range("Q11").Value = 601573
$ws.Range("R11").Value = 6992600
$ws.Range("X11").Value = "2025_0868"
$ws.Range("Z11").Value = "13:18"
$ws.Range("AB11").Value = "13:18"
$ws.Range("AX11").Value = "David Isaksson"

# Row 12 <- source row 11 (id 131106319)
$ws.Range("A12").Value = 131106319
$ws.Range("B12").Value = 92107
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 658
$ws.Range("F12").Value = "Rosenticka"
$ws.Range("G12").Value = "Fomitopsis rosea"
$ws.Range("H12").Value = "(Alb. & Schwein.:Fr.) P.Karst."
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = "2"
$ws.Range("J12").Value = "mycel"
$ws.Range("Q12").Value = 601569
$ws.Range("R12").Value = 6992657
$ws.Range("X12").Value = "2025_0864"
$ws.Range("Z12").Value = "13:14"
$ws.Range("AB12").Value = "13:14"
$ws.Range("AX12").Value = "Alexander Hoffmann"

# Row 13 <- source row 14 (id 131106325)
$ws.Range("A13").Value = 131106325
$ws.Range("B13").Value = 91809
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 1202
$ws.Range("F13").Value = "Ullticka"
$ws.Range("G13").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H13").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("I13").Value = ""
$ws.Range("J13").Value = "mycel"
$ws.Range("Q13").Value = 601615
$ws.Range("R13").Value = 6992785
$ws.Range("X13").Value = "2025_0858"
$ws.Range("Z13").Value = "12:21"
$ws.Range("AB13").Value = "12:21"
$ws.Range("AX13").Value = "Alexander Hoffmann"

# Row 14 <- source row 13 (id 131106312)
$ws.Range("A14").Value = 131106312
$ws.Range("B14").Value = 92107
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 658
$ws.Range("F14").Value = "Rosenticka"
$ws.Range("G14").Value = "Fomitopsis rosea"
$ws.Range("H14").Value = "(Alb. & Schwein.:Fr.) P.Karst."
$ws.Range("I14").Value = ""
$ws.Range("J14").Value = "mycel"
$ws.Range("Q14").Value = 601540
$ws.Range("R14").Value = 6992576
$ws.Range("X14").Value = "2025_0872"
$ws.Range("Z14").Value = "13:29"
$ws.Range("AB14").Value = "13:29"
$ws.Range("AX14").Value = "Alexander Hoffmann"

# Row 23 <- source row 25 (id 131106323)
$ws.Range("A23").Value = 131106323
$ws.Range("B23").Value = 92107
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 658
$ws.Range("F23").Value = "Rosenticka"
$ws.Range("G23").Value = "Fomitopsis rosea"
$ws.Range("H23").Value = "(Alb. & Schwein.:Fr.) P.Karst."
$ws.Range("I23").Value = ""
$ws.Range("J23").Value = ""
$ws.Range("Q23").Value = 601607
$ws.Range("R23").Value = 6992738
$ws.Range("X23").Value = "2025_0860"
$ws.Range("Z23").Value = "12:35"
$ws.Range("AB23").Value = "12:35"
$ws.Range("AX23").Value = "David Isaksson"

# Row 24 <- source row 23 (id 131106327)
$ws.Range("A24").Value = 131106327
$ws.Range("B24").Value = 91809
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 1202
$ws.Range("F24").Value = "Ullticka"
$ws.Range("G24").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H24").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("I24").Value = ""
$ws.Range("J24").Value = "mycel"
$ws.Range("Q24").Value = 601607
$ws.Range("R24").Value = 6992789
$ws.Range("X24").Value = "2025_0856"
$ws.Range("Z24").Value = "12:10"
$ws.Range("AB24").Value = "12:10"
$ws.Range("AX24").Value = "David Isaksson, Alexander Hoffmann"

# Row 25 <- source row 24 (id 131106311)
$ws.Range("A25").Value = 131106311
$ws.Range("B25").Value = 91809
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 1202
$ws.Range("F25").Value = "Ullticka"
$ws.Range("G25").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H25").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("I25").Value = ""
$ws.Range("J25").Value = ""
$ws.Range("Q25").Value = 601498
$ws.Range("R25").Value = 6992583
$ws.Range("X25").Value = "2025_0873"
$ws.Range("Z25").Value = "13:32"
$ws.Range("AB25").Value = "13:32"
$ws.Range("AX25").Value = "Alexander Hoffmann"

